$d = $word.ActiveDocument

# --- Header block text replacements (unique strings, safe for Find/Replace) ---
$d.Content.Find.Execute("Client name: Bruce Banner", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Client name: John Doe", 2)

# NOTE: this replacement contains a straight apostrophe ('). Word's Find/Replace
# "replace with" text passes through AutoFormat/AutoCorrect smart-quote handling,
# which would silently turn it into a curly quote. Locate the range with Find
# (no replacement text) and then assign Range.Text directly so the literal
# straight apostrophe is preserved verbatim.
$rngTherapist = $d.Content
$foundTherapist = $rngTherapist.Find.Execute("Therapist providing services: Jacquie Atkins, LPC")
if ($foundTherapist) {
    $rngTherapist.Text = "Therapist providing services: Ryan O'Farrell, Psy.D."
}

$d.Content.Find.Execute("Date of Estimate: 08/01/2022", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Date of Estimate: 08/02/2022", 2)

# --- Table updates (use explicit Cell addressing to avoid ambiguous duplicate values) ---
# Table 1: 12-session itemized estimate
$t1 = $d.Tables.Item(1)
$t1.Cell(3, 4).Range.Text = "165"      # Initial evaluation - Cost
$t1.Cell(3, 6).Range.Text = "165"      # Initial evaluation - Estimate
$t1.Cell(4, 2).Range.Text = "90837"    # Psychotherapy - Service code
$t1.Cell(4, 4).Range.Text = "165"      # Psychotherapy - Cost
$t1.Cell(4, 6).Range.Text = "1980"     # Psychotherapy - Estimate

# Table 2: 24-session itemized estimate
$t2 = $d.Tables.Item(2)
$t2.Cell(3, 4).Range.Text = "165"      # Initial evaluation - Cost
$t2.Cell(3, 6).Range.Text = "165"      # Initial evaluation - Estimate
$t2.Cell(4, 2).Range.Text = "90837"    # Psychotherapy - Service code
$t2.Cell(4, 4).Range.Text = "165"      # Psychotherapy - Cost
$t2.Cell(4, 6).Range.Text = "3960"     # Psychotherapy - Estimate

# --- Estimate range summary text ---
$d.Content.Find.Execute("Estimate range: `$1845-3525", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "Estimate range: `$2170-4150", 2)
